# The sheet contains a header row (row 1) followed by 30 data rows
# (rows 2-31, columns A-T). The edit reorders these 30 data rows
# (the set of rows is identical, only their order changes) without
# touching the header row or any styling.
#
# Mapping below: for each destination data-row index (0-based, where
# 0 = sheet row 2), which source data-row index (0-based) supplies the
# values.
$map = @(22,23,24,25,0,16,5,6,7,8,14,17,1,2,26,27,18,19,20,21,9,10,11,12,15,3,4,28,29,13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 31
$numCols = 20

# 1) Snapshot every current cell value for the data rows into memory
#    before we overwrite anything.
$snapshot = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value())
    }
    $snapshot += ,$rowVals
}

# 2) Write the snapshot back out in the new order.
for ($i = 0; $i -lt $map.Count; $i++) {
    $srcIdx = $map[$i]
    $destRow = $firstDataRow + $i
    $srcRowVals = $snapshot[$srcIdx]
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcRowVals[$c - 1]
    }
}
